$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.845000000000001
$ws.Range("B4").Value = 6.318
$ws.Range("C6").Value = -12.52
$ws.Range("B7").Value = 6.411
$ws.Range("C7").Value = -12.918
$ws.Range("B8").Value = 5.823
$ws.Range("C8").Value = -12.241
$ws.Range("A11").Value = -21.643
$ws.Range("E11").Value = 12.763
$ws.Range("A12").Value = -21.456
$ws.Range("B12").Value = 6.695
$ws.Range("B14").Value = 7.545
$ws.Range("E14").Value = 12.924
$ws.Range("A15").Value = -21.096
$ws.Range("C19").Value = -12.235
$ws.Range("E19").Value = 12.929
$ws.Range("C21").Value = -12.725
$ws.Range("E21").Value = 12.919
$ws.Range("B22").Value = 6.724000000000001
$ws.Range("C24").Value = -12.255
$ws.Range("C25").Value = -12.69

$wb.Save()
